$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Diem tong ket" (final score) for group N10 (row 11):
# B11 changes from "88.89/400" to "400/400"
$ws.Range("B11").Value = "400/400"

# D11 changes from 4 to 10
$ws.Range("D11").Value = 10

# Move the active selection to D11 (matches the saved selection state in the file)
$ws.Range("D11").Select()
